$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Fix title typo: "Requirments" -> "Requirements"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Requirments", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Requirements", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Merge the three runs describing "Full Name" into a single run (the
#    rendered text is unchanged, but Find/Replace across the run boundaries
#    naturally collapses them into one run, matching the target OOXML).
# ---------------------------------------------------------------------------
$oldFullName = " this should be First Name and Last Name combined, and should not be stored  in the database"
$d.Content.Find.Execute($oldFullName, $true, $false, $false, $false, $false,
                         $true, 1, $false, $oldFullName, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Move the word "Project " from the start of the "5. Update the ..."
#    paragraph up onto the end of the preceding "Tests.Unit.Sample.Services"
#    heading paragraph, and give the now-"5. Update the ..." paragraph the
#    (new) "No Spacing" paragraph style.
# ---------------------------------------------------------------------------
$text = $d.Content.Text
$anchor = $text.IndexOf("Tests.Unit.Sample.Services")
$paraMark = $text.IndexOf([char]13, $anchor)
$d.Range($paraMark, $paraMark).InsertBefore("Project ") | Out-Null

$text = $d.Content.Text
$moved = $text.IndexOf("Project 5. Update the")
$d.Range($moved, $moved + 8).Text = ""

$text = $d.Content.Text
$target = $text.IndexOf("5. Update the EmployeeServiceTest")
$d.Range($target, $target).Paragraphs(1).Range.Style = "No Spacing"

$noSpacing = $d.Styles("No Spacing")
$noSpacing.Priority = 1
$noSpacing.ParagraphFormat.SpaceAfter = 0
$noSpacing.ParagraphFormat.LineSpacingRule = 0

# ---------------------------------------------------------------------------
# 4. Move the word "Project" (no trailing space) from the start of the
#    "6. Update the ..." paragraph up onto the end of the preceding
#    "Tests.Unit.Sample.Repositories" heading paragraph.
# ---------------------------------------------------------------------------
$text = $d.Content.Text
$anchor2 = $text.IndexOf("Tests.Unit.Sample.Repositories")
$paraMark2 = $text.IndexOf([char]13, $anchor2)
$d.Range($paraMark2, $paraMark2).InsertBefore("Project") | Out-Null

$text = $d.Content.Text
$moved2 = $text.IndexOf("Project 6. Update the")
$d.Range($moved2, $moved2 + 7).Text = ""

# ---------------------------------------------------------------------------
# 5. "6). Web API Project 7. " -> "6). Web API Project. "
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("6). Web API Project 7. ", $true, $false, $false,
                         $false, $false, $true, 1, $false,
                         "6). Web API Project. ", 2) | Out-Null
